$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the two runs " " and "XAMPP is regularly updated to the
# latest releases of " (in the "Virtualised server" section) into a single
# run, leaving everything around them untouched.
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute(" XAMPP is regularly updated to the latest releases of ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c1Start = $find1.Start
$c1End = $find1.End

# Temporary bookmarks act as hard boundaries so that the upcoming
# find/replace only re-merges the runs strictly inside the match (and does
# not bleed into neighbouring runs that happen to share the same, empty,
# run formatting).
$d.Bookmarks.Add("ZZ_Barrier1L", $d.Range($c1Start, $c1Start)) | Out-Null
$d.Bookmarks.Add("ZZ_Barrier1R", $d.Range($c1End, $c1End)) | Out-Null

$find1b = $d.Content
$find1b.Find.Execute(" XAMPP is regularly updated to the latest releases of ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " XAMPP is regularly updated to the latest releases of ", 2)

$d.Bookmarks("ZZ_Barrier1L").Delete()
$d.Bookmarks("ZZ_Barrier1R").Delete()

# ---------------------------------------------------------------------------
# Change 2: merge the two runs "XAMPP is offered in both a full and a
# standard version (Smaller version)" and "." into a single run, keeping
# "another. " before it and " We are using ..." after it untouched.
# ---------------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("XAMPP is offered in both a full and a standard version (Smaller version)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c2Start = $find2.Start
$c2End = $find2.End
$c2PeriodEnd = $c2End + 1

$d.Bookmarks.Add("ZZ_Barrier2L", $d.Range($c2Start, $c2Start)) | Out-Null
$d.Bookmarks.Add("ZZ_Barrier2R", $d.Range($c2PeriodEnd, $c2PeriodEnd)) | Out-Null

$find2b = $d.Content
$find2b.Find.Execute("XAMPP is offered in both a full and a standard version (Smaller version).", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "XAMPP is offered in both a full and a standard version (Smaller version).", 2)

$d.Bookmarks("ZZ_Barrier2L").Delete()
$d.Bookmarks("ZZ_Barrier2R").Delete()

# ---------------------------------------------------------------------------
# Change 3: in the "Image editing" heading line, insert a new run containing
# a single space between the existing " –" run and the _GoBack bookmark that
# precedes "Microsoft Photos", so the text reads " – Microsoft Photos".
# ---------------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute("Image editing –", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$c3InsertPos = $find3.End

$ins = $d.Range($c3InsertPos, $c3InsertPos)
$ins.InsertAfter(" ")

# Give the freshly inserted space a momentary explicit format and then clear
# it again; this stops the engine from silently re-absorbing the new run
# back into the preceding " –" run, while leaving no visible formatting on
# the final text.
$newRun = $d.Range($c3InsertPos, $c3InsertPos + 1)
$newRun.Bold = 1
$newRun2 = $d.Range($c3InsertPos, $c3InsertPos + 1)
$newRun2.Bold = 0
